$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'87.780.51"
$ws.Range("E2").Value = "'  -2.49%  "
$ws.Range("D2:E2").Style = "Normal"

$ws.Range("D3").Value = "'3.057.61"
$ws.Range("E3").Value = "'  -5.09%  "
$ws.Range("D3:E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'208.63"
$ws.Range("E5").Value = "'  -4.52%  "
$ws.Range("D5:E5").Style = "Normal"

$ws.Range("D6").Value = "'615.59"
$ws.Range("E6").Value = "'  -2.70%  "
$ws.Range("D6:E6").Style = "Normal"

$ws.Range("D7").Value = "'0.365"
$ws.Range("E7").Value = "'  -7.46%  "
$ws.Range("D7:E7").Style = "Normal"

$ws.Range("D8").Value = "'0.806"
$ws.Range("E8").Value = "'  +15.64%  "
$ws.Range("D8:E8").Style = "Normal"

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "'  +0.02%  "
$ws.Range("D9:E9").Style = "Normal"

$ws.Range("D10").Value = "'3.053.90"
$ws.Range("E10").Value = "'  -5.11%  "
$ws.Range("D10:E10").Style = "Normal"

$ws.Range("D11").Value = "'0.596"
$ws.Range("E11").Value = "'  +3.44%  "
$ws.Range("D11:E11").Style = "Normal"

$ws.Range("E12").Value = "'  -1.00%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.0000233"
$ws.Range("E13").Value = "'  -11.08%  "
$ws.Range("D13:E13").Style = "Normal"

$ws.Range("D14").Value = "'5.23"
$ws.Range("E14").Value = "'  -3.70%  "
$ws.Range("D14:E14").Style = "Normal"

$ws.Range("D15").Value = "'87.707.02"
$ws.Range("E15").Value = "'  -2.29%  "
$ws.Range("D15:E15").Style = "Normal"

$ws.Range("D16").Value = "'3.615.59"
$ws.Range("E16").Value = "'  -5.21%  "
$ws.Range("D16:E16").Style = "Normal"

$ws.Range("D17").Value = "'31.46"
$ws.Range("E17").Value = "'  -6.76%  "
$ws.Range("D17:E17").Style = "Normal"

$ws.Range("D18").Value = "'3.034.70"
$ws.Range("E18").Value = "'  -5.79%  "
$ws.Range("D18:E18").Style = "Normal"

$ws.Range("D19").Value = "'3.17"
$ws.Range("E19").Value = "'  -9.21%  "
$ws.Range("D19:E19").Style = "Normal"

$ws.Range("D20").Value = "'13.07"
$ws.Range("E20").Value = "'  -4.11%  "
$ws.Range("D20:E20").Style = "Normal"

$ws.Range("D21").Value = "'0.0000195"
$ws.Range("E21").Value = "'  -20.26%  "
$ws.Range("D21:E21").Style = "Normal"

$ws.Range("D22").Value = "'416.71"
$ws.Range("E22").Value = "'  -5.79%  "
$ws.Range("D22:E22").Style = "Normal"

$ws.Range("D23").Value = "'8.05"
$ws.Range("E23").Value = "'  -7.48%  "
$ws.Range("D23:E23").Style = "Normal"

$ws.Range("D24").Value = "'4.82"
$ws.Range("E24").Value = "'  -5.87%  "
$ws.Range("D24:E24").Style = "Normal"

$ws.Range("D25").Value = "'5.32"
$ws.Range("E25").Value = "'  +0.76%  "
$ws.Range("D25:E25").Style = "Normal"

$ws.Range("D26").Value = "'11.51"
$ws.Range("E26").Value = "'  -4.60%  "
$ws.Range("D26:E26").Style = "Normal"

$ws.Range("D27").Value = "'80.87"
$ws.Range("E27").Value = "'  -2.31%  "
$ws.Range("D27:E27").Style = "Normal"

$ws.Range("D28").Value = "'3.214.83"
$ws.Range("E28").Value = "'  -5.31%  "
$ws.Range("D28:E28").Style = "Normal"

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "'  -0.12%  "
$ws.Range("D29:E29").Style = "Normal"

$ws.Range("D30").Value = "'1.08"
$ws.Range("E30").Value = "'  +8.59%  "
$ws.Range("D30:E30").Style = "Normal"

$ws.Range("D31").Value = "'0.162"
$ws.Range("E31").Value = "'  +1.16%  "
$ws.Range("D31:E31").Style = "Normal"

$ws.Range("D32").Value = "'7.98"
$ws.Range("E32").Value = "'  -7.41%  "
$ws.Range("D32:E32").Style = "Normal"

$ws.Range("D33").Value = "'496.87"
$ws.Range("E33").Value = "'  -9.37%  "
$ws.Range("D33:E33").Style = "Normal"

$ws.Range("D34").Value = "'3.48"
$ws.Range("E34").Value = "'  -17.64%  "
$ws.Range("D34:E34").Style = "Normal"

$ws.Range("D35").Value = "'6.53"
$ws.Range("E35").Value = "'  -7.84%  "
$ws.Range("D35:E35").Style = "Normal"

$ws.Range("E36").Value = "'  -7.52%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.22"
$ws.Range("E37").Value = "'  -7.63%  "
$ws.Range("D37:E37").Style = "Normal"

$ws.Range("D38").Value = "'21.96"
$ws.Range("E38").Value = "'  -2.34%  "
$ws.Range("D38:E38").Style = "Normal"

$ws.Range("E39").Value = "'  -0.44%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  -0.83%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +0.29%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  -0.02%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.356"
$ws.Range("E43").Value = "'  -5.81%  "
$ws.Range("D43:E43").Style = "Normal"

$ws.Range("D44").Value = "'146.97"
$ws.Range("E44").Value = "'  -0.08%  "
$ws.Range("D44:E44").Style = "Normal"

$ws.Range("E45").Value = "'  -8.99%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.132"
$ws.Range("E46").Value = "'  +5.32%  "
$ws.Range("D46:E46").Style = "Normal"

$ws.Range("D47").Value = "'43.34"
$ws.Range("E47").Value = "'  -1.35%  "
$ws.Range("D47:E47").Style = "Normal"

$ws.Range("D48").Value = "'0.0658"
$ws.Range("E48").Value = "'  +7.26%  "
$ws.Range("D48:E48").Style = "Normal"

$ws.Range("D49").Value = "'154.14"
$ws.Range("E49").Value = "'  -11.64%  "
$ws.Range("D49:E49").Style = "Normal"

$ws.Range("D50").Value = "'0.693"
$ws.Range("E50").Value = "'  -8.67%  "
$ws.Range("D50:E50").Style = "Normal"

$ws.Range("E51").Value = "'  -8.60%  "
$ws.Range("E51").Style = "Normal"
